$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 21.04553260289762
$ws.Range("B3").Value = 22.0070710644361
$ws.Range("B4").Value = 65.31612035052051
$ws.Range("B5").Value = 44.17008018759426
$ws.Range("B6").Value = 47.492330679556

$ws.Range("A7").Value = "bus"
$ws.Range("B7").Value = 75.3056429174495
$ws.Range("A8").Value = "bus"
$ws.Range("B8").Value = 43.1133569440285
$ws.Range("A9").Value = "bus"
$ws.Range("B9").Value = 41.47088942768687
$ws.Range("A10").Value = "bus"
$ws.Range("B10").Value = 8.16923076923082
$ws.Range("A11").Value = "bus"
$ws.Range("B11").Value = 58.04980978904219

$ws.Range("A12").Value = "truck"
$ws.Range("B12").Value = 62.04096999573512
$ws.Range("A13").Value = "truck"
$ws.Range("B13").Value = 52.65503583354378
$ws.Range("A14").Value = "truck"
$ws.Range("B14").Value = 59.56242979952989
$ws.Range("A15").Value = "truck"
$ws.Range("B15").Value = 56.17233960436053

$ws.Range("B16").Value = 16.61538461538419

$ws.Range("A17").Value = "motorcycle"
$ws.Range("B17").Value = 82.98241558561205

$ws.Range("B18").Value = 56.88216042630249
$ws.Range("B19").Value = 24.15930719924168

$ws.Range("A20").Value = "motorcycle"
$ws.Range("B20").Value = 33.33350842526671
$ws.Range("A21").Value = "motorcycle"
$ws.Range("B21").Value = 63.72674657089972
